# Applies the TDY cash-flow update:
#  - narrows columns C and D to match the width used by the rest of the
#    quarterly columns (was inherited from column B, now ~14.3 chars)
#  - fills in the previously-blank "B" column (most recent quarter) for
#    every line item that didn't have it yet
#  - corrects several mis-keyed historical figures in columns C-F

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: columns C and D should match the narrower width
#     already used by column E onward (was same as column B). ColumnWidth
#     is expressed in characters; Excel stores/quantizes it in Maximum
#     Digit Width (MDW) pixel units, so this is the closest achievable
#     character-width setting to the target 14.3 stored width.
$ws.Columns.Item(3).ColumnWidth = 13.6
$ws.Columns.Item(4).ColumnWidth = 13.6

# --- Fill in previously-empty column B (most recent quarter) values ---
$ws.Range("B3").Value = 116200000.0
$ws.Range("B4").Value = 15900000.0
$ws.Range("B5").Value = 45200000.0
$ws.Range("B7").Value = -21000000.0
$ws.Range("B10").Value = -10200000.0
$ws.Range("B11").Value = 400000.0
$ws.Range("B12").Value = 667400000.0
$ws.Range("B13").Value = -68800000.0
$ws.Range("B16").Value = -67900000.0
$ws.Range("B17").Value = 2351000000.0
$ws.Range("B18").Value = 36900000.0
$ws.Range("B19").Value = 14800000.0
$ws.Range("B20").Value = 2382300000.0
$ws.Range("B21").Value = 21000000.0
$ws.Range("B22").Value = 3002800000.0
$ws.Range("B23").Value = 673100000.0
$ws.Range("B24").Value = 3675900000.0
$ws.Range("B25").Value = 27400000.0
$ws.Range("B27").Value = 130500000.0
$ws.Range("B28").Value = 36900000.0
$ws.Range("B29").Value = 36900000.0

# --- Correct mis-keyed historical figures ---
$ws.Range("F4").Value = -15200000.0

$ws.Range("B6").Value = 398700000.0
$ws.Range("C6").Value = 54300000.0
$ws.Range("D6").Value = 42000000.0
$ws.Range("E6").Value = 16300000.0
$ws.Range("F6").Value = 8200000.0

$ws.Range("B8").Value = 941000000.0
$ws.Range("C8").Value = 954000000.0
$ws.Range("D8").Value = 734800000.0
$ws.Range("E8").Value = 540500000.0
$ws.Range("F8").Value = 295900000.0

$ws.Range("F27").Value = -29200000.0
